$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header labels on row 2: the previously "unnamed" placeholder
# headers for the totals columns should read "total".
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
